# "error solve ifrs list"
# Rewrites the per-year IFRS-consolidated figures in the company_list sheet
# (rows 2-9, columns D:AJ) with corrected values, and removes the stray
# trailing metric columns / rows that no longer apply.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row target data: cells to overwrite with a corrected value ("Set"),
# and cells that must become empty ("Clear") because the corrected figures
# no longer populate them.
$rowData = @{
    2 = @{ Set = @{ "D"=256; "E"=-74; "F"=-74; "G"=-227; "H"=-227; "I"=-223; "J"=-5; "K"=462; "L"=298; "M"=165; "N"=159; "O"=6; "P"=74; "Q"=-47; "R"=-5; "S"=71; "T"=63; "U"=-109; "V"=275; "W"=-29.07; "X"=-88.97; "AA"=180.6; "AB"=259.25; "AC"=-743; "AD"=-1.68; "AE"=428; "AF"=2.91; "AG"=0; "AH"=0; "AI"=0; "AJ"=37097139 }; Clear = @("Y", "Z") }
    3 = @{ Set = @{ "D"=398; "E"=-47; "F"=-92; "G"=-28; "H"=-25; "I"=-19; "J"=-6; "K"=789; "L"=412; "M"=377; "N"=377; "O"=0; "P"=128; "Q"=-40; "R"=30; "S"=16; "T"=13; "U"=-53; "V"=359; "W"=-11.76; "X"=-6.33; "Y"=-7.17; "Z"=-4.03; "AA"=109.13; "AB"=278.19; "AC"=-39; "AD"=-24.65; "AE"=589; "AF"=1.63; "AG"=0; "AH"=0; "AI"=0; "AJ"=64101882 }; Clear = @() }
    4 = @{ Set = @{ "D"=360; "E"=-28; "F"=-28; "G"=-110; "H"=-141; "I"=-139; "J"=-2; "K"=387; "L"=104; "M"=283; "N"=279; "O"=4; "P"=145; "Q"=-23; "R"=85; "S"=-18; "T"=5; "U"=-29; "V"=85; "W"=-7.74; "X"=-39.07; "Y"=-42.25; "Z"=-23.95; "AA"=36.54; "AB"=166.45; "AC"=-206; "AD"=-3.45; "AE"=384; "AF"=1.85; "AG"=0; "AH"=0; "AI"=0; "AJ"=72690174 }; Clear = @() }
    5 = @{ Set = @{ "D"=844; "E"=-45; "F"=-45; "G"=7; "H"=7; "I"=7; "J"=0; "K"=998; "L"=308; "M"=691; "N"=686; "O"=4; "P"=245; "Q"=-197; "R"=-459; "S"=682; "T"=0; "U"=-198; "V"=290; "W"=-5.33; "X"=0.84; "Y"=1.49; "Z"=1.02; "AA"=44.52; "AB"=258.19; "AC"=6; "AD"=57.76; "AE"=560; "AF"=0.64; "AG"=0; "AH"=0; "AI"=0; "AJ"=122690174 }; Clear = @() }
    6 = @{ Set = @{ "D"=1213; "E"=-34; "F"=-34; "G"=-67; "H"=-67; "I"=-67; "K"=800; "L"=61; "M"=739; "N"=734; "P"=484; "Q"=-76; "R"=123; "S"=-102; "T"=10; "U"=-86; "V"=47; "W"=-2.84; "X"=-5.53; "Y"=-9.44; "Z"=-7.46; "AA"=8.300000000000001; "AB"=99.76000000000001; "AC"=-36; "AD"=-5.41; "AE"=303; "AF"=0.64; "AG"=0; "AH"=0; "AI"=0; "AJ"=242188696 }; Clear = @("J", "O") }
    7 = @{ Set = @{}; Clear = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ") }
    8 = @{ Set = @{}; Clear = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ") }
    9 = @{ Set = @{}; Clear = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ") }
}

foreach ($r in $rowData.Keys) {
    $entry = $rowData[$r]

    foreach ($col in $entry.Set.Keys) {
        $ws.Range("$col$r").Value = $entry.Set[$col]
    }

    foreach ($col in $entry.Clear) {
        $ws.Range("$col$r").ClearContents()
    }
}
